# Scheduled runner update: refresh cached Market Board pricing/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) on a handful of rows across
# several sheets. Values with no computable profit have their profit cell(s)
# cleared instead of holding a stale number.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H9").Value = 106.36842
$ws.Range("I9").Value = 63.9375
$ws.Range("J9").Value = 332.66666
$ws.Range("K9").Value = 63.9375
$ws.Range("L9").Value = 332.66666
$ws.Range("M9").Value = 105.0625
$ws.Range("N9").Value = -670.66666

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H105").Value = 38000
$ws.Range("J105").Value = 38000
$ws.Range("L105").Value = 38000
$ws.Range("N105").Value = -44988

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H110").Value = 42000
$ws.Range("J110").Value = 42000
$ws.Range("L110").Value = 42000
$ws.Range("N110").Value = -50180

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H120").Value = 46000
$ws.Range("J120").Value = 46000
$ws.Range("L120").Value = 46000
$ws.Range("N120").Value = -55676

$ws.Range("H132").Value = 2098437.2
$ws.Range("I132").Value = 1702.6046
$ws.Range("J132").Value = 11114397
$ws.Range("K132").Value = 5107.8138
$ws.Range("L132").Value = 33343191
$ws.Range("M132").Value = -2577.8138
$ws.Range("N132").Value = -33348251

$ws.Range("H135").Value = 17857816
$ws.Range("I135").Value = 647.1923
$ws.Range("J135").Value = 250001000
$ws.Range("K135").Value = 5824.7307
$ws.Range("L135").Value = 2250009000
$ws.Range("M135").Value = -3289.7307
$ws.Range("N135").Value = -2250014070

$ws.Range("H138").Value = 3811.25
$ws.Range("I138").Value = 3879.7334
$ws.Range("J138").Value = 3775.8276
$ws.Range("K138").Value = 11639.2002
$ws.Range("L138").Value = 11327.4828
$ws.Range("M138").Value = -6499.200199999999
$ws.Range("N138").Value = -21607.4828

$ws.Range("H141").Value = 1462.8868
$ws.Range("I141").Value = 847.61224
$ws.Range("K141").Value = 2542.83672
$ws.Range("M141").Value = 2637.16328

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 620.9
$ws.Range("I97").Value = 618.1667
$ws.Range("J97").Value = 625
$ws.Range("K97").Value = 618.1667
$ws.Range("L97").Value = 625
$ws.Range("M97").Value = -122.1667
$ws.Range("N97").Value = -1617

$ws.Range("H132").Value = 19161870
$ws.Range("J132").Value = 9262159
$ws.Range("L132").Value = 27786477
$ws.Range("N132").Value = -27791537

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1668133.4
$ws.Range("I107").Value = 2001660
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 2001660
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = -1999740
$ws.Range("N107").Value = -4340

$ws.Range("H134").Value = 23292948
$ws.Range("I134").Value = 27778844
$ws.Range("K134").Value = 83336532
$ws.Range("M134").Value = -83333997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1692783.2
$ws.Range("I31").Value = 1191.0769
$ws.Range("K31").Value = 1191.0769
$ws.Range("M31").Value = -896.0769

$ws.Range("H34").Value = 1692783.2
$ws.Range("I34").Value = 1191.0769
$ws.Range("K34").Value = 1191.0769
$ws.Range("M34").Value = -989.0769

$ws.Range("H58").Value = 1467376.9
$ws.Range("I58").Value = 1123.1111
$ws.Range("K58").Value = 1123.1111
$ws.Range("M58").Value = -920.1111000000001

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H132").Value = 881.4286
$ws.Range("I132").Value = 693.9636
$ws.Range("K132").Value = 2081.8908
$ws.Range("M132").Value = 448.1091999999999

$ws.Range("H134").Value = 910167.1
$ws.Range("I134").Value = 1081.325
$ws.Range("J134").Value = 10001025
$ws.Range("K134").Value = 3243.975
$ws.Range("L134").Value = 30003075
$ws.Range("M134").Value = -708.9750000000004
$ws.Range("N134").Value = -30008145

$ws.Range("H136").Value = 1467376.9
$ws.Range("I136").Value = 1123.1111
$ws.Range("K136").Value = 3369.3333
$ws.Range("M136").Value = -819.3333000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 230.92857
$ws.Range("I26").Value = 59.391304
$ws.Range("J26").Value = 1020
$ws.Range("K26").Value = 178.173912
$ws.Range("L26").Value = 3060
$ws.Range("M26").Value = 109.826088
$ws.Range("N26").Value = -3636

$ws.Range("H107").Value = 950443.4
$ws.Range("J107").Value = 1162
$ws.Range("L107").Value = 3486
$ws.Range("N107").Value = -7326

$ws.Range("H122").Value = 1186.7097
$ws.Range("I122").Value = 286.46155
$ws.Range("J122").Value = 1836.8889
$ws.Range("K122").Value = 2578.15395
$ws.Range("L122").Value = 16532.0001
$ws.Range("M122").Value = -128.1539499999999
$ws.Range("N122").Value = -21432.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 755
$ws.Range("I41").Value = 755
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 755
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -400
$ws.Range("N41").ClearContents()

$ws.Range("H132").Value = 9957668
$ws.Range("I132").Value = 9905351
$ws.Range("J132").Value = 10102993
$ws.Range("K132").Value = 29716053
$ws.Range("L132").Value = 30308979
$ws.Range("M132").Value = -29713523
$ws.Range("N132").Value = -30314039

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3664525
$ws.Range("I132").Value = 4202631.5
$ws.Range("K132").Value = 12607894.5
$ws.Range("M132").Value = -12605364.5

$ws.Range("H136").Value = 2925100.5
$ws.Range("I136").Value = 3004130.2
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 9012390.600000001
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -9009840.600000001
$ws.Range("N136").Value = -8100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 32000
$ws.Range("J95").Value = 32000
$ws.Range("L95").Value = 32000
$ws.Range("N95").Value = -37492
